# WAT New Test script added WAT157,158,159
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Copy formatting from row 124 (style pattern A/B/C/E = border only, D = border+fill)
# so the new rows 126-128 visually match the rest of the table.
$ws.Range("A124:E124").Copy()
$ws.Range("A126:E126").PasteSpecial(-4122)
$ws.Range("A124:E124").Copy()
$ws.Range("A127:E127").PasteSpecial(-4122)
$ws.Range("A124:E124").Copy()
$ws.Range("A128:E128").PasteSpecial(-4122)

# Row 126 - WAT157 / WAT-1210
$ws.Cells.Item(126, 1).Value = "WAT157"
$ws.Cells.Item(126, 2).Value = "WAT-1210"
$ws.Cells.Item(126, 3).Value = "Verify that ‘Select All’ option should not display when quantity of search results of an author morethan 50."
$ws.Cells.Item(126, 4).Value = "Y"

# Row 127 - WAT158 / WAT-1211
$ws.Cells.Item(127, 1).Value = "WAT158"
$ws.Cells.Item(127, 2).Value = "WAT-1211"
$ws.Cells.Item(127, 3).Value = "Verify that ‘Select All’ option should display when quantity of search results of an author lessthan 50."
$ws.Cells.Item(127, 4).Value = "Y"

# Row 128 - WAT159 / WAT-1333
$ws.Cells.Item(128, 1).Value = "WAT159"
$ws.Cells.Item(128, 2).Value = "WAT-1333"
$ws.Cells.Item(128, 3).Value = "Verify that ‘Select All’ option should not display when quantity of search results of an author Only 1"
$ws.Cells.Item(128, 4).Value = "Y"

# Match the recorded selection/active cell from the edit
[void]$ws.Range("A126:E128").Select()
